# Word COM-interop script implementing the target diff:
#  1. Collapse the three CORE COMPETENCIES detail paragraphs into a single
#     condensed summary paragraph.
#  2. Append a new "TECHNICAL SKILLS" section (heading + three detail
#     paragraphs) at the end of the document, just before the sectPr.

$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- 1. Condense CORE COMPETENCIES paragraphs -----------------------------
# Paragraph 6 ("Product Management & Strategy: ...") is rewritten in place
# with the short summary text; paragraphs 7 ("Technical Product
# Development: ...") and 8 ("Platform & Infrastructure: ...") are then
# deleted outright (their detail, in condensed form, moves to the new
# TECHNICAL SKILLS section appended below).

$coreCompParagraph = $d.Paragraphs.Item(6)
$coreCompParagraph.Range.Text = "Product Management & Strategy " + $bullet + " Technical Product Development " + $bullet + " Platform & Infrastructure"

$d = $word.ActiveDocument
$d.Paragraphs.Item(7).Range.Delete()
$d = $word.ActiveDocument
$d.Paragraphs.Item(7).Range.Delete()

# --- 2. Append new TECHNICAL SKILLS section -------------------------------
# Helper pattern: collapse to the end of the document, insert a bare new
# paragraph mark, explicitly set the new (still-empty) paragraph's style
# BEFORE filling in its text, so the new paragraph doesn't inherit
# whatever style the previous last paragraph had.

function Add-TrailingParagraph([string]$styleName, [string]$text) {
    $doc = $word.ActiveDocument
    $r = $doc.Paragraphs.Last.Range
    $r.Collapse(0)
    $r.InsertAfter("`r")
    $doc = $word.ActiveDocument
    $p = $doc.Paragraphs.Last
    $p.Style = $styleName
    $p.Range.Text = $text
}

Add-TrailingParagraph "Heading2" "TECHNICAL SKILLS"
Add-TrailingParagraph "Normal" "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics"
Add-TrailingParagraph "Normal" "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration"
Add-TrailingParagraph "Normal" "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training"

$d = $word.ActiveDocument
Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
